$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 26, shifting rows 26-30 down to 27-31
$ws.Rows.Item(26).Insert()

# Copy formatting of the date cell in new row 26 from the row below (old row 26, now row 27)
$ws.Range("D26").NumberFormat = $ws.Range("D27").NumberFormat

# Populate the new row 26 with data
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 45021
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 100112043
$ws.Range("G26").Value = "Pepino dulce"
$ws.Range("H26").Value = "Cultivar IV Región"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 270
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17500
$ws.Range("N26").Value = "`$/bandeja 18 kilos"
$ws.Range("O26").Value = "Provincia de Limarí"
$ws.Range("P26").Value = 972
$ws.Range("Q26").Value = 18
$ws.Range("R26").Value = "Hortaliza"
